$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
  @{Row=2; D="29.949.08"; E="  -1.58%  "},
  @{Row=3; D="1.893.51"; E="  -2.40%  "},
  @{Row=4; D="0.9994"; E="  -0.62%  "},
  @{Row=5; D="0.7371"; E="  -1.52%  "},
  @{Row=6; D="242.92"; E="  -1.26%  "},
  @{Row=7; D="0.9991"; E="  -0.59%  "},
  @{Row=8; D="0.3106"},
  @{Row=9; D="26.29"},
  @{Row=10; D="0.06909"; E="  -1.12%  "},
  @{Row=11; D="0.7731"; E="  -1.11%  "},
  @{Row=12; D="0.07956"; E="  -0.48%  "},
  @{Row=13; D="1.901.09"; E="  -1.97%  "},
  @{Row=14; D="5.228"; E="  -2.43%  "},
  @{Row=15; D="91.57"; E="  -3.19%  "},
  @{Row=16; E="  -1.96%  "},
  @{Row=17; D="29.956.75"; E="  -1.55%  "},
  @{Row=18; D="5.774"; E="  +0.12%  "},
  @{Row=19; D="240.61"; E="  -5.07%  "},
  @{Row=20; E="  -2.15%  "},
  @{Row=21; D="0.9993"; E="  -0.43%  "},
  @{Row=22; D="2.124.23"; E="  -3.07%  "},
  @{Row=23; D="0.9993"; E="  -0.59%  "},
  @{Row=24; D="6.907"; E="  +3.35%  "},
  @{Row=25; D="9.310"; E="  -2.26%  "},
  @{Row=26; D="164.73"; E="  -0.33%  "},
  @{Row=27; D="18.87"; E="  -0.82%  "},
  @{Row=28; D="0.1272"; E="  -3.78%  "},
  @{Row=29; D="2.017"; E="  -10.59%  "},
  @{Row=30; D="1.362"; E="  -0.23%  "},
  @{Row=31; D="1.532"; E="  +1.05%  "},
  @{Row=32; D="4.312"; E="  -1.40%  "},
  @{Row=33; D="4.063"},
  @{Row=34; D="0.05112"; E="  -1.01%  "},
  @{Row=35; D="1.282"; E="  +0.22%  "},
  @{Row=36; D="0.7370"; E="  -1.27%  "},
  @{Row=37; D="2.714"; E="  -2.44%  "},
  @{Row=38; D="0.01919"; E="  -1.77%  "},
  @{Row=39; D="2.776"; E="  -1.08%  "},
  @{Row=40; D="6.313"; E="  -1.80%  "},
  @{Row=41; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="74.36"; E="  -5.30%  "},
  @{Row=42; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.4461"; E="  -0.32%  "},
  @{Row=43; D="1.936"; E="  -1.69%  "},
  @{Row=44; D="0.9997"; E="  -0.52%  "},
  @{Row=45; D="0.8395"; E="  +0.99%  "},
  @{Row=46; D="7.639"; E="  +2.09%  "},
  @{Row=47; D="100.95"; E="  -0.33%  "},
  @{Row=48; D="9.810"; E="  +0.21%  "},
  @{Row=49; D="36.59"; E="  -1.75%  "},
  @{Row=50; D="2.021.23"; E="  -3.12%  "},
  @{Row=51; D="939.30"; E="  -4.05%  "}
)

foreach ($item in $changes) {
  $r = $item.Row
  if ($item.ContainsKey("B")) { $ws.Range("B$r").Value = $item.B }
  if ($item.ContainsKey("C")) { $ws.Range("C$r").Value = $item.C }
  if ($item.ContainsKey("D")) {
    $ws.Range("D$r").Value = "'" + $item.D
    $ws.Range("D$r").ClearFormats()
  }
  if ($item.ContainsKey("E")) { $ws.Range("E$r").Value = $item.E }
}
